$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 with the same formatting (style) as H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I and J for rows 2-11
$iValues = @(4, 8, 4, 8, 5, 7, 6, 4, 5, 7)
$jValues = @(6, 9, 5, 8, 5, 7, 6, 4, 5, 7)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
